$wb = $excel.ActiveWorkbook

# Row 28 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7414.4287
$ws.Range("I28").Value = 7973.5386
$ws.Range("K28").Value = 7973.5386
$ws.Range("M28").Value = -7488.5386

# Row 41 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 257.8
$ws.Range("I41").Value = 324.66666
$ws.Range("J41").Value = 157.5
$ws.Range("K41").Value = 324.66666
$ws.Range("L41").Value = 157.5
$ws.Range("M41").Value = 115.33334
$ws.Range("N41").Value = -1037.5

# Row 43 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 44241.668
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 97043.75
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 97043.75
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -97181.75

# Row 74 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5900
$ws.Range("I74").Value = 5900
$ws.Range("K74").Value = 5900
$ws.Range("M74").Value = -4964

# Row 77 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5900
$ws.Range("I77").Value = 5900
$ws.Range("K77").Value = 29500
$ws.Range("M77").Value = -24820

# Row 111 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 703.2222
$ws.Range("I111").Value = 704.2857
$ws.Range("J111").Value = 699.5
$ws.Range("K111").Value = 2112.8571
$ws.Range("L111").Value = 2098.5
$ws.Range("M111").Value = 954.1428999999998
$ws.Range("N111").Value = -8232.5

# Row 138 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1685.4286
$ws.Range("I138").Value = 899
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 2697
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = 2443
$ws.Range("N138").Value = -16280

# Row 2 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1619.1765
$ws.Range("I2").Value = 1569.1333
$ws.Range("J2").Value = 1994.5
$ws.Range("K2").Value = 1569.1333
$ws.Range("L2").Value = 1994.5
$ws.Range("M2").Value = -1456.1333
$ws.Range("N2").Value = -2220.5

# Row 110 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2497
$ws.Range("I110").Value = 655.6667
$ws.Range("J110").Value = 3878
$ws.Range("K110").Value = 655.6667
$ws.Range("L110").Value = 3878
$ws.Range("M110").Value = 1389.3333
$ws.Range("N110").Value = -7968

# Row 116 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1619.1765
$ws.Range("I116").Value = 1569.1333
$ws.Range("J116").Value = 1994.5
$ws.Range("K116").Value = 1569.1333
$ws.Range("L116").Value = 1994.5
$ws.Range("M116").Value = 724.8667
$ws.Range("N116").Value = -6582.5

# Row 3 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1619.1765
$ws.Range("I3").Value = 1569.1333
$ws.Range("J3").Value = 1994.5
$ws.Range("K3").Value = 1569.1333
$ws.Range("L3").Value = 1994.5
$ws.Range("M3").Value = -1455.1333
$ws.Range("N3").Value = -2222.5

# Row 134 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7616.4736
$ws.Range("I134").Value = 7571.353
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 22714.059
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -20179.059
$ws.Range("N134").Value = -29070

# Row 10 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 278.2857
$ws.Range("I10").Value = 158
$ws.Range("K10").Value = 158
$ws.Range("M10").Value = -19

# Row 19 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 577.38464
$ws.Range("I19").Value = 139.55556
$ws.Range("J19").Value = 1562.5
$ws.Range("K19").Value = 139.55556
$ws.Range("L19").Value = 1562.5
$ws.Range("M19").Value = 30.44443999999999
$ws.Range("N19").Value = -1902.5

# Row 24 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 577.38464
$ws.Range("I24").Value = 139.55556
$ws.Range("J24").Value = 1562.5
$ws.Range("K24").Value = 139.55556
$ws.Range("L24").Value = 1562.5
$ws.Range("M24").Value = 30.44443999999999
$ws.Range("N24").Value = -1902.5

# Row 127 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 10490
$ws.Range("J127").Value = 10490
$ws.Range("L127").Value = 31470
$ws.Range("N127").Value = -41390

# Row 131 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1302.3
$ws.Range("I131").Value = 1005.6
$ws.Range("J131").Value = 1599
$ws.Range("K131").Value = 3016.8
$ws.Range("L131").Value = 4797
$ws.Range("M131").Value = 2023.2
$ws.Range("N131").Value = -14877

# Row 140 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2978.3076
$ws.Range("I140").Value = 2610.818
$ws.Range("K140").Value = 7832.454000000001
$ws.Range("M140").Value = -2652.454000000001

# Row 43 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 24220.875
$ws.Range("J43").Value = 29294.5
$ws.Range("L43").Value = 29294.5
$ws.Range("N43").Value = -29596.5

# Row 46 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 20000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 70 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4833
$ws.Range("I70").Value = 2250
$ws.Range("K70").Value = 2250
$ws.Range("M70").Value = -1980

# Row 73 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4833
$ws.Range("I73").Value = 2250
$ws.Range("K73").Value = 2250
$ws.Range("M73").Value = -1314

# Row 22 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1256.9
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1341
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 1341
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1931

# Row 27 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1256.9
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1341
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 1341
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -1555

# Row 55 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 197.625
$ws.Range("J55").Value = 149.2
$ws.Range("L55").Value = 149.2
$ws.Range("N55").Value = -495.2

# Row 100 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1525
$ws.Range("I100").Value = 1033.3334
$ws.Range("K100").Value = 1033.3334
$ws.Range("M100").Value = -492.3334

# Row 136 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3503.3333
$ws.Range("I136").Value = 3503.3333
$ws.Range("K136").Value = 10509.9999
$ws.Range("M136").Value = -7959.999899999999

# Row 6 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 125
$ws.Range("J6").Value = 125
$ws.Range("L6").Value = 125
$ws.Range("N6").Value = -355

# Row 31 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 30000
$ws.Range("J31").Value = 30000
$ws.Range("K31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = -29652
$ws.Range("N31").Value = -30696

# Row 122 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1127.4375
$ws.Range("I122").Value = 1127.4375
$ws.Range("K122").Value = 3382.3125
$ws.Range("M122").Value = -932.3125

# Row 126 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2927.875
$ws.Range("I126").Value = 2862.4167
$ws.Range("J126").Value = 3124.25
$ws.Range("K126").Value = 8587.250100000001
$ws.Range("L126").Value = 9372.75
$ws.Range("M126").Value = -6117.250100000001
$ws.Range("N126").Value = -14312.75

# Row 136 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3355.1667
$ws.Range("I136").Value = 3427.4
$ws.Range("J136").Value = 2994
$ws.Range("K136").Value = 10282.2
$ws.Range("L136").Value = 8982
$ws.Range("M136").Value = -7732.200000000001
$ws.Range("N136").Value = -14082
